$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.016973376274109
$ws.Range("B1").Value = 1.357475280761719
$ws.Range("C1").Value = 2.207844495773315
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.970597267150879
